$d = $word.ActiveDocument

$replacements = @(
    @{old="910×8=7280"; new="306×9=2754"},
    @{old="390×5=1950"; new="516×5=2580"},
    @{old="338×3=1014"; new="142×2=284"},
    @{old="349×8=2792"; new="601×5=3005"},
    @{old="826×4=3304"; new="835×2=1670"},
    @{old="477×3=1431"; new="529×8=4232"},
    @{old="109×8=872"; new="305×6=1830"},
    @{old="120×9=1080"; new="235×4=940"},
    @{old="488×5=2440"; new="620×7=4340"},
    @{old="176×8=1408"; new="238×9=2142"},
    @{old="208×4=832"; new="838×3=2514"},
    @{old="241×7=1687"; new="834×8=6672"},
    @{old="194×9=1746"; new="723×9=6507"},
    @{old="194×2=388"; new="845×3=2535"},
    @{old="805×2=1610"; new="338×8=2704"},
    @{old="825×8=6600"; new="231×3=693"},
    @{old="163×9=1467"; new="279×3=837"},
    @{old="570×5=2850"; new="906×5=4530"},
    @{old="334×9=3006"; new="483×6=2898"},
    @{old="419×2=838"; new="878×6=5268"},
    @{old="497×2=994"; new="849×7=5943"},
    @{old="478×5=2390"; new="868×7=6076"},
    @{old="386×8=3088"; new="507×4=2028"},
    @{old="686×5=3430"; new="574×3=1722"},
    @{old="965×3=2895"; new="955×8=7640"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
